$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 49
$ws.Range("H3").Value = 320
$ws.Range("H4").Value = 132
$ws.Range("H5").Value = 585
$ws.Range("H6").Value = 79
$ws.Range("H7").Value = 96
$ws.Range("H8").Value = 92
$ws.Range("H10").Value = 215
$ws.Range("H11").Value = 72
$ws.Range("H12").Value = 266
$ws.Range("H14").Value = 66
$ws.Range("H15").Value = 468
$ws.Range("H16").Value = 103
